$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 8, duplicating the content of row 5 (same course entry re-listed)
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "Durga Raju"
$ws.Range("C8").Value = "Mastering SQL using Postgresql"
$ws.Range("D8").Value = "https://www.udemy.com/course/mastering-sql-using-postgresql/?couponCode=KEEPLEARNING"
$ws.Range("E8").Value = "Video Course"

# Update the selection / view to reflect where the user left the cursor
$ws.Range("C6").Select()
$excel.ActiveWindow.ScrollColumn = 3
